$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 0.277
$ws.Range("E2").Value = 0.0493
$ws.Range("I2").Value = -1.611009174311927
$ws.Range("J2").Value = -1.239425874622378
$ws.Range("K2").Value = 5.08
$ws.Range("L2").Value = 9.321100917431192
$ws.Range("M2").Value = 0.0
$ws.Range("N2").Value = 0.0
$ws.Range("O2").Value = 0.0
$ws.Range("P2").Value = 0.0
$ws.Range("Q2").Value = 0.0
$ws.Range("R2").Value = 0.0
$ws.Range("U2").Value = 0.915
$ws.Range("V2").Value = 0.03009868421052632
$ws.Range("W2").Value = 0.05893271461716938
$ws.Range("X2").Value = 0.03756833553445338
$ws.Range("Y2").Value = 0.021364379082716
$ws.Range("Z2").Value = 0.006446576217456619
$ws.Range("AA2").Value = -0.007990053366640989
$ws.Range("AB2").Value = 0.03755771546865593
$ws.Range("AC2").Value = -0.04554776883529692
$ws.Range("AD2").Value = 0.029
$ws.Range("AF2").Value = 0.029
$ws.Range("AG2").Value = -0.886
$ws.Range("AH2").Value = 0.0009530382201189656
$ws.Range("AI2").Value = 0.0003414616915305726
$ws.Range("AJ2").Value = -0.03001965169072305
$ws.Range("AK2").Value = -0.0105458614040517
$ws.Range("AL2").Value = 0.003
$ws.Range("AM2").Value = -0.011
$ws.Range("AO2").Value = -292.6666666666667
$ws.Range("AQ2").Value = 79.81818181818183

# Row 3 updates
$ws.Range("D3").Value = 0.277
$ws.Range("E3").Value = 0.0493
$ws.Range("I3").Value = -1.611009174311927
$ws.Range("J3").Value = -1.239425874622378
$ws.Range("K3").Value = 5.08
$ws.Range("L3").Value = 9.321100917431192
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("U3").Value = 0.915
$ws.Range("V3").Value = 0.03009868421052632
$ws.Range("W3").Value = 0.05893271461716938
$ws.Range("X3").Value = 0.03756833553445338
$ws.Range("Y3").Value = 0.021364379082716
$ws.Range("Z3").Value = 0.006446576217456619
$ws.Range("AA3").Value = -0.007990053366640989
$ws.Range("AB3").Value = 0.03755771546865593
$ws.Range("AC3").Value = -0.04554776883529692
$ws.Range("AD3").Value = 0.029
$ws.Range("AF3").Value = 0.029
$ws.Range("AG3").Value = -0.886
$ws.Range("AH3").Value = 0.0009530382201189656
$ws.Range("AI3").Value = 0.0003414616915305726
$ws.Range("AJ3").Value = -0.03001965169072305
$ws.Range("AK3").Value = -0.0105458614040517
$ws.Range("AL3").Value = 0.003
$ws.Range("AM3").Value = -0.011
$ws.Range("AO3").Value = -292.6666666666667
$ws.Range("AQ3").Value = 79.81818181818183

# Remove T2 / T3 cells entirely (column dropped)
$ws.Range("T2").ClearContents()
$ws.Range("T3").ClearContents()

Write-Host "Applied capital structure database updates"
